$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "SceneQuestId2" column header (K3) to "SceneQuestId" -
# part of wiring up the basic task/quest system's scene-quest link column.
$ws.Range("K3").Value = "SceneQuestId"

# Fill in the first quest row's ResetOnLeave flag (J4) and its new
# SceneQuestId value (K4). J5 already stores the word "true" as literal
# text (the column is text-formatted), so copy that cell's value across
# instead of assigning the literal string, which would otherwise get
# auto-coerced to a boolean True.
$ws.Range("J5").Copy()
$ws.Range("J4").PasteSpecial(-4163)
$ws.Range("K4").Value = 42120017

# Move the active selection to J4, matching where the author left off.
[void]$ws.Range("J4").Select()
